# Weekly refresh of the "Ají" (Hortaliza, Mapocho Venta Directa de Santiago)
# price sheet: the data rows (2,4-18; row 3 and row 10 stay put) get
# reshuffled to a new order. Implemented as a row permutation: snapshot every
# source row's values first (so overlapping reads/writes don't clobber each
# other), then write each destination row from its snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new_row -> old_row (where the data for the new row comes from)
$mapping = @{
    2  = 6
    4  = 2
    5  = 17
    6  = 18
    7  = 9
    8  = 7
    9  = 13
    11 = 14
    12 = 15
    13 = 4
    14 = 11
    15 = 5
    16 = 12
    17 = 16
    18 = 8
}

# Snapshot the current contents of every row referenced above before any
# writes happen, since several rows are both a source and a destination.
$snapshot = @{}
foreach ($oldRow in $mapping.Values) {
    if (-not $snapshot.ContainsKey($oldRow)) {
        $snapshot[$oldRow] = $ws.Range("A$oldRow`:R$oldRow").Value2
    }
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $ws.Range("A$newRow`:R$newRow").Value2 = $snapshot[$oldRow]
}
